$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studenti24_genova")

# New rows appended at the bottom of the manual (row 53 intentionally left blank)
$ws.Range("A54").Value = "Variabili da ricodificare recode"
$ws.Range("A55").Value = "comp_sost_2"
$ws.Range("A56").Value = "comp_sost_6"
$ws.Range("A57").Value = "collab_classe2_2"
$ws.Range("A58").Value = "collab_classe2_5"

# Column widths (inputs picked so the engine's char->pixel rounding lands on
# the closest representable width to the authored 35.5546875 / 152)
$ws.Columns.Item(1).ColumnWidth = 34.666666666666664
$ws.Columns.Item(2).ColumnWidth = 151.16666666666666

# View / selection state
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("A59").Select()
